$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Init_unrecReturn")

# Update header info cells
$ws.Range("B2").Value = "C6:D9"
$ws.Range("B3").Value = "all"
$ws.Range("B4").Value = 20200630

# Update the year / DeferredReturn table (shift years forward, set DeferredReturn flags to 1)
$ws.Range("C7").Value = 2020
$ws.Range("D7").Value = 1

$ws.Range("C8").Value = 2021
$ws.Range("D8").Value = 1

$ws.Range("C9").Value = 2022
$ws.Range("D9").Value = 1

# Row 10 loses its year value (C10 cleared) but D10 stays blank
$ws.Range("C10").ClearContents()
$ws.Range("D10").ClearContents()

# Row 11 is removed entirely - delete it so rows shift up
$ws.Rows.Item(11).Delete()

# Update selection to B3
$ws.Range("B3").Select()
